{"js": "// Update the worksheet's date heading and the 25 division-fact answers\n// laid out in a 5x5 grid that lives on every 4th row (0, 4, 8, 12, 16)\n// of the single table in the document body.\n\n// 1) Date heading: first paragraph of the body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.getRange().insertText(\"2024-01-21 Sunday\", \"Replace\");\n\n// 2) Division-fact grid values, addressed by (row, column) so the edit is\n//    robust even though some new values duplicate other OLD values\n//    elsewhere in the table (e.g. new row4/col2 \"75\u00f78=9, 3\" equals the old\n//    row0/col0 value).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"97\u00f79=10, 7\", \"19\u00f76=3, 1\", \"46\u00f79=5, 1\", \"35\u00f75=7, 0\", \"75\u00f79=8, 3\"],\n  4: [\"61\u00f75=12, 1\", \"25\u00f73=8, 1\", \"75\u00f78=9, 3\", \"91\u00f76=15, 1\", \"23\u00f72=11, 1\"],\n  8: [\"82\u00f73=27, 1\", \"68\u00f72=34, 0\", \"96\u00f75=19, 1\", \"36\u00f78=4, 4\", \"45\u00f78=5, 5\"],\n  12: [\"64\u00f78=8, 0\", \"49\u00f77=7, 0\", \"21\u00f73=7, 0\", \"64\u00f72=32, 0\", \"13\u00f78=1, 5\"],\n  16: [\"43\u00f74=10, 3\", \"50\u00f72=25, 0\", \"87\u00f75=17, 2\", \"94\u00f79=10, 4\", \"77\u00f76=12, 5\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const row = Number(rowIndex);\n  const rowValues = newValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(row, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet's date heading and the 25 division-fact answers\n# laid out in a 5x5 grid that lives on every 4th row (1, 5, 9, 13, 17 in\n# 1-based COM indexing) of the single table in the document body.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading: first paragraph of the body.\n$d.Paragraphs.Item(1).Range.Text = \"2024-01-21 Sunday\"\n\n# 2) Division-fact grid values, addressed by (row, column) so the edit is\n#    robust even though some new values duplicate other OLD values\n#    elsewhere in the table (e.g. new row5/col3 \"75\u00f78=9, 3\" equals the old\n#    row1/col1 value).\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"97\u00f79=10, 7\"\n$t.Cell(1, 2).Range.Text = \"19\u00f76=3, 1\"\n$t.Cell(1, 3).Range.Text = \"46\u00f79=5, 1\"\n$t.Cell(1, 4).Range.Text = \"35\u00f75=7, 0\"\n$t.Cell(1, 5).Range.Text = \"75\u00f79=8, 3\"\n\n$t.Cell(5, 1).Range.Text = \"61\u00f75=12, 1\"\n$t.Cell(5, 2).Range.Text = \"25\u00f73=8, 1\"\n$t.Cell(5, 3).Range.Text = \"75\u00f78=9, 3\"\n$t.Cell(5, 4).Range.Text = \"91\u00f76=15, 1\"\n$t.Cell(5, 5).Range.Text = \"23\u00f72=11, 1\"\n\n$t.Cell(9, 1).Range.Text = \"82\u00f73=27, 1\"\n$t.Cell(9, 2).Range.Text = \"68\u00f72=34, 0\"\n$t.Cell(9, 3).Range.Text = \"96\u00f75=19, 1\"\n$t.Cell(9, 4).Range.Text = \"36\u00f78=4, 4\"\n$t.Cell(9, 5).Range.Text = \"45\u00f78=5, 5\"\n\n$t.Cell(13, 1).Range.Text = \"64\u00f78=8, 0\"\n$t.Cell(13, 2).Range.Text = \"49\u00f77=7, 0\"\n$t.Cell(13, 3).Range.Text = \"21\u00f73=7, 0\"\n$t.Cell(13, 4).Range.Text = \"64\u00f72=32, 0\"\n$t.Cell(13, 5).Range.Text = \"13\u00f78=1, 5\"\n\n$t.Cell(17, 1).Range.Text = \"43\u00f74=10, 3\"\n$t.Cell(17, 2).Range.Text = \"50\u00f72=25, 0\"\n$t.Cell(17, 3).Range.Text = \"87\u00f75=17, 2\"\n$t.Cell(17, 4).Range.Text = \"94\u00f79=10, 4\"\n$t.Cell(17, 5).Range.Text = \"77\u00f76=12, 5\"\n"}
